# Natmi following Dr Hou advice
# Update Calr-Scarf1 LR-pair statistics with recalculated values (3 expressing cells
# instead of 1) for Ligand-expressing cells / Receptor-expressing cells and all
# downstream derived expression / specificity metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3, 1, 93.32574466666667, 279.977234, 0.2327963689879921, 0.2327963689879922, 3, 1, 10.16187433333333, 30.485623, 0.6867885401361519, 0.686788540136152, 948.3644893674202, 8535.280404306783, 0.15988187840626, 0.1598818784062601)
    3  = @(3, 1, 93.32574466666667, 279.977234, 0.2327963689879921, 0.2327963689879922, 3, 1, 0.5923493333333333, 1.777048, 0.0400338284597913, 0.04003382845979131, 55.28144263613688, 497.532983725232, 0.009319729902127556, 0.009319729902127559)
    4  = @(3, 1, 93.32574466666667, 279.977234, 0.2327963689879921, 0.2327963689879922, 3, 1, 2.456299666666667, 7.368899, 0.1660085931857371, 0.1660085931857371, 229.2359955161518, 2063.123959645366, 0.03864619771444432, 0.03864619771444434)
    5  = @(3, 1, 93.32574466666667, 279.977234, 0.2327963689879921, 0.2327963689879922, 3, 1, 1.585696666666667, 4.75709, 0.1071690382183197, 0.1071690382183197, 147.9863222321178, 1331.87690008906, 0.02494856296516018, 0.0249485629651602)
    6  = @(3, 1, 164.7897643333334, 494.369293, 0.4110597662007076, 0.4110597662007077, 3, 1, 10.16187433333333, 30.485623, 0.6867885401361519, 0.686788540136152, 1674.572876574949, 15071.15588917454, 0.2823111367376919, 0.282311136737692)
    7  = @(3, 1, 164.7897643333334, 494.369293, 0.4110597662007076, 0.4110597662007077, 3, 1, 0.5923493333333333, 1.777048, 0.0400338284597913, 0.04003382845979131, 97.61310704300712, 878.517963387064, 0.01645629616680105, 0.01645629616680105)
    8  = @(3, 1, 164.7897643333334, 494.369293, 0.4110597662007076, 0.4110597662007077, 3, 1, 2.456299666666667, 7.368899, 0.1660085931857371, 0.1660085931857371, 404.7730432020453, 3642.957388818407, 0.06823945350223747, 0.06823945350223748)
    9  = @(3, 1, 164.7897643333334, 494.369293, 0.4110597662007076, 0.4110597662007077, 3, 1, 1.585696666666667, 4.75709, 0.1071690382183197, 0.1071690382183197, 261.3065800041522, 2351.75922003737, 0.0440528797939772, 0.04405287979397721)
    10 = @(3, 1, 95.835818, 287.507454, 0.2390576204784642, 0.2390576204784643, 3, 1, 10.16187433333333, 30.485623, 0.6867885401361519, 0.686788540136152, 973.8715391482048, 8764.843852333843, 0.1641820341768267, 0.1641820341768268)
    11 = @(3, 1, 95.835818, 287.507454, 0.2390576204784642, 0.2390576204784643, 3, 1, 0.5923493333333333, 1.777048, 0.0400338284597913, 0.04003382845979131, 56.76828290175467, 510.914546115792, 0.00957039177024073, 0.009570391770240733)
    12 = @(3, 1, 95.835818, 287.507454, 0.2390576204784642, 0.2390576204784643, 3, 1, 2.456299666666667, 7.368899, 0.1660085931857371, 0.1660085931857371, 235.4014878081273, 2118.613390273146, 0.0396856192659597, 0.03968561926595971)
    13 = @(3, 1, 95.835818, 287.507454, 0.2390576204784642, 0.2390576204784643, 3, 1, 1.585696666666667, 4.75709, 0.1071690382183197, 0.1071690382183197, 151.9665371498733, 1367.69883434886, 0.0256195752654371, 0.02561957526543711)
    14 = @(3, 1, 46.93870866666666, 140.816126, 0.117086244332836, 0.117086244332836, 3, 1, 10.16187433333333, 30.485623, 0.6867885401361519, 0.686788540136152, 476.9852588396109, 4292.867329556498, 0.08041349081537319, 0.08041349081537323)
    15 = @(3, 1, 46.93870866666666, 140.816126, 0.117086244332836, 0.117086244332836, 3, 1, 0.5923493333333333, 1.777048, 0.0400338284597913, 0.04003382845979131, 27.80411278622755, 250.237015076048, 0.004687410620621967, 0.004687410620621968)
    16 = @(3, 1, 46.93870866666666, 140.816126, 0.117086244332836, 0.117086244332836, 3, 1, 2.456299666666667, 7.368899, 0.1660085931857371, 0.1660085931857371, 115.2955344516971, 1037.659810065274, 0.01943732270309558, 0.01943732270309559)
    17 = @(3, 1, 46.93870866666666, 140.816126, 0.117086244332836, 0.117086244332836, 3, 1, 1.585696666666667, 4.75709, 0.1071690382183197, 0.1071690382183197, 74.4305538703711, 669.8749848333399, 0.01254802019374522, 0.01254802019374522)
}

foreach ($r in $data.Keys) {
    $values = $data[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 5 + $i   # column E is index 5
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}
